# Add a new service event entry to the "Card15" sheet.
# - Fills the previously-blank placeholder cells of row 18 with the literal
#   text "nan" (matching the sheet's convention for empty data points).
# - Appends a brand-new row 19 describing the newly logged service event.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card15")

# --- Row 18: materialize the blank cells as literal "nan" text ---------
$nanCols = @("B","C","D","E","F","G","H","I","J","K","M","P","Q","R")
foreach ($col in $nanCols) {
    $ws.Range($col + "18").Value = "nan"
}

# --- Row 19: new service entry -----------------------------------------
# First stamp out a fully-populated, still-blank row so every column from
# A to R has a real (empty) cell, matching the layout used by every other
# row in this sheet.
$ws.Range("A1000:R1000").Copy($ws.Range("A19:R19"))

# "15" must be stored as text (like the rest of column A), so force the
# cell to Text format before writing the value; otherwise Excel would
# auto-convert the numeric-looking string into a number.
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "15"

$ws.Range("L19").Value = "20\10\2025"
$ws.Range("N19").Value = "تم تغيير الجرائد الخلفيه (1_5_8)"
$ws.Range("O19").Value = "الخبير"
